$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the formatting (date style) from row 31 column A down to the new rows 32:43
$ws.Range("A31").Copy()
$ws.Range("A32:A43").PasteSpecial(-4122)

# Write the updated timestamp / import / export values for rows 2 through 43
$ws.Cells.Item(2, 1).Value = 46022
$ws.Cells.Item(2, 2).Value = 3.738
$ws.Cells.Item(2, 3).Value = 0.352
$ws.Cells.Item(3, 1).Value = 46022.01041666666
$ws.Cells.Item(3, 2).Value = 0.169
$ws.Cells.Item(3, 3).Value = 1.205
$ws.Cells.Item(4, 1).Value = 46022.02083333334
$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(4, 3).Value = 11.834
$ws.Cells.Item(5, 1).Value = 46022.03125
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = 15.21
$ws.Cells.Item(6, 1).Value = 46022.04166666666
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = 18.679
$ws.Cells.Item(7, 1).Value = 46022.05208333334
$ws.Cells.Item(7, 2).Value = 1.855
$ws.Cells.Item(7, 3).Value = 1.725
$ws.Cells.Item(8, 1).Value = 46022.0625
$ws.Cells.Item(8, 2).Value = 0.773
$ws.Cells.Item(8, 3).Value = 3.118
$ws.Cells.Item(9, 1).Value = 46022.07291666666
$ws.Cells.Item(9, 2).Value = 3.014
$ws.Cells.Item(9, 3).Value = 0.007
$ws.Cells.Item(10, 1).Value = 46022.08333333334
$ws.Cells.Item(10, 2).Value = 4.141
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(11, 1).Value = 46022.09375
$ws.Cells.Item(11, 2).Value = 21.322
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(12, 1).Value = 46022.10416666666
$ws.Cells.Item(12, 2).Value = 43.196
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(13, 1).Value = 46022.11458333334
$ws.Cells.Item(13, 2).Value = 42.58
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(14, 1).Value = 46022.125
$ws.Cells.Item(14, 2).Value = 13.042
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(15, 1).Value = 46022.13541666666
$ws.Cells.Item(15, 2).Value = 36.942
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(16, 1).Value = 46022.14583333334
$ws.Cells.Item(16, 2).Value = 19.712
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(17, 1).Value = 46022.15625
$ws.Cells.Item(17, 2).Value = 3.38
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(18, 1).Value = 46022.16666666666
$ws.Cells.Item(18, 2).Value = 1.307
$ws.Cells.Item(18, 3).Value = 2.068
$ws.Cells.Item(19, 1).Value = 46022.17708333334
$ws.Cells.Item(19, 2).Value = 0
$ws.Cells.Item(19, 3).Value = 26.829
$ws.Cells.Item(20, 1).Value = 46022.1875
$ws.Cells.Item(20, 2).Value = 0
$ws.Cells.Item(20, 3).Value = 47.755
$ws.Cells.Item(21, 1).Value = 46022.19791666666
$ws.Cells.Item(21, 2).Value = 0
$ws.Cells.Item(21, 3).Value = 23.409
$ws.Cells.Item(22, 1).Value = 46022.20833333334
$ws.Cells.Item(22, 2).Value = 0
$ws.Cells.Item(22, 3).Value = 24.673
$ws.Cells.Item(23, 1).Value = 46022.21875
$ws.Cells.Item(23, 2).Value = 0
$ws.Cells.Item(23, 3).Value = 13.163
$ws.Cells.Item(24, 1).Value = 46022.22916666666
$ws.Cells.Item(24, 2).Value = 3.6
$ws.Cells.Item(24, 3).Value = 0.416
$ws.Cells.Item(25, 1).Value = 46022.23958333334
$ws.Cells.Item(25, 2).Value = 15.688
$ws.Cells.Item(25, 3).Value = 0.02
$ws.Cells.Item(26, 1).Value = 46022.25
$ws.Cells.Item(26, 2).Value = 0
$ws.Cells.Item(26, 3).Value = 24.243
$ws.Cells.Item(27, 1).Value = 46022.26041666666
$ws.Cells.Item(27, 2).Value = 0
$ws.Cells.Item(27, 3).Value = 37.087
$ws.Cells.Item(28, 1).Value = 46022.27083333334
$ws.Cells.Item(28, 2).Value = 0
$ws.Cells.Item(28, 3).Value = 18.589
$ws.Cells.Item(29, 1).Value = 46022.28125
$ws.Cells.Item(29, 2).Value = 0
$ws.Cells.Item(29, 3).Value = 25.015
$ws.Cells.Item(30, 1).Value = 46022.29166666666
$ws.Cells.Item(30, 2).Value = 0
$ws.Cells.Item(30, 3).Value = 18.029
$ws.Cells.Item(31, 1).Value = 46022.30208333334
$ws.Cells.Item(31, 2).Value = 0.004
$ws.Cells.Item(31, 3).Value = 4.444
$ws.Cells.Item(32, 1).Value = 46022.3125
$ws.Cells.Item(32, 2).Value = 0.008
$ws.Cells.Item(32, 3).Value = 4.168
$ws.Cells.Item(33, 1).Value = 46022.32291666666
$ws.Cells.Item(33, 2).Value = 0.012
$ws.Cells.Item(33, 3).Value = 10.778
$ws.Cells.Item(34, 1).Value = 46022.33333333334
$ws.Cells.Item(34, 2).Value = 0.291
$ws.Cells.Item(34, 3).Value = 9.536
$ws.Cells.Item(35, 1).Value = 46022.34375
$ws.Cells.Item(35, 2).Value = 0.05
$ws.Cells.Item(35, 3).Value = 6.59
$ws.Cells.Item(36, 1).Value = 46022.35416666666
$ws.Cells.Item(36, 2).Value = 0.092
$ws.Cells.Item(36, 3).Value = 4.598
$ws.Cells.Item(37, 1).Value = 46022.36458333334
$ws.Cells.Item(37, 2).Value = 0
$ws.Cells.Item(37, 3).Value = 13.382
$ws.Cells.Item(38, 1).Value = 46022.375
$ws.Cells.Item(38, 2).Value = 6.584
$ws.Cells.Item(38, 3).Value = 0.163
$ws.Cells.Item(39, 1).Value = 46022.38541666666
$ws.Cells.Item(39, 2).Value = 0.719
$ws.Cells.Item(39, 3).Value = 7.978
$ws.Cells.Item(40, 1).Value = 46022.39583333334
$ws.Cells.Item(40, 2).Value = 0
$ws.Cells.Item(40, 3).Value = 26.509
$ws.Cells.Item(41, 1).Value = 46022.40625
$ws.Cells.Item(41, 2).Value = 0
$ws.Cells.Item(41, 3).Value = 40.748
$ws.Cells.Item(42, 1).Value = 46022.41666666666
$ws.Cells.Item(42, 2).Value = 5.782
$ws.Cells.Item(42, 3).Value = 0.97
$ws.Cells.Item(43, 1).Value = 46022.42708333334
$ws.Cells.Item(43, 2).Value = 0
$ws.Cells.Item(43, 3).Value = 0
